$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update date (D2), Volumen (M2), Precio minimo/maximo/promedio (N2:P2),
# Unidad de comercializacion (Q2), Precio $/Kg (S2), Kg/unidad (T2)
$ws.Range("D2").Value2 = 44973
$ws.Range("M2").Value2 = 60
$ws.Range("N2").Value2 = 12000
$ws.Range("O2").Value2 = 12000
$ws.Range("P2").Value2 = 12000
$ws.Range("Q2").Value = "`$/bandeja 5 kilos"
$ws.Range("S2").Value2 = 2400
$ws.Range("T2").Value2 = 5

# Row 3: update date (D3), Volumen (M3), Precio minimo/maximo/promedio (N3:P3),
# Unidad de comercializacion (Q3), Precio $/Kg (S3), Kg/unidad (T3)
$ws.Range("D3").Value2 = 44238
$ws.Range("M3").Value2 = 35
$ws.Range("N3").Value2 = 20000
$ws.Range("O3").Value2 = 20000
$ws.Range("P3").Value2 = 20000
$ws.Range("Q3").Value = "`$/bandeja 10 kilos"
$ws.Range("S3").Value2 = 2000
$ws.Range("T3").Value2 = 10
